# A new weekly price record was added to the "Ajo" (garlic) price log for
# Vega Monumental Concepción. The new observation is inserted as row 13,
# pushing all existing rows (13-113) down by one (to 14-114).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing row 13 and below shift to 14+
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new price observation.
$ws.Cells.Item(13, 1).Value  = 11
$ws.Cells.Item(13, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value  = "Bíobío"
$ws.Cells.Item(13, 4).Value  = 44537
$ws.Cells.Item(13, 5).Value  = 8
$ws.Cells.Item(13, 6).Value  = 100112003
$ws.Cells.Item(13, 7).Value  = "Ajo"
$ws.Cells.Item(13, 8).Value  = "Chino"
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 310
$ws.Cells.Item(13, 11).Value = 17000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 17516
$ws.Cells.Item(13, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(13, 15).Value = "China"
$ws.Cells.Item(13, 16).Value = 1752
$ws.Cells.Item(13, 17).Value = 10
$ws.Cells.Item(13, 18).Value = "Hortaliza"
